$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the Saturday header to Sunday (shared string text update)
$ws.Range("H2").Value = "Domingos  1/2"

# Change the fill color of that header cell to the theme "Accent 2" color
$ws.Range("H2").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorAccent2

# Fill in the missing "Abs" entry for that day
$ws.Range("H4").Value = "Abs"

# Widen column H slightly so the new label fits
$ws.Columns("H").ColumnWidth = 12.022135416666666

# Update the active cell selection to match
$ws.Range("H14").Select() | Out-Null
